# Atualização de bases das ligas, do dia: 30-05-2024 às 12:21
# Swap the data (columns B:AD) between pairs of rows. Column A (the
# sequential row id) stays put; all other fields move together as a
# full record swap between the two matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two parallel flat arrays describing the row pairs to swap (nested
# arrays such as @(@(19,20), @(25,26)) get flattened when enumerated in
# this environment, so parallel arrays + an index loop are used instead).
$rowsA = @(19, 25, 42, 107, 130, 141, 157, 161)
$rowsB = @(20, 26, 43, 108, 131, 142, 158, 163)

for ($i = 0; $i -lt $rowsA.Count; $i++) {
    $r1 = $rowsA[$i]
    $r2 = $rowsB[$i]

    $range1 = $ws.Range("B${r1}:AD${r1}")
    $range2 = $ws.Range("B${r2}:AD${r2}")

    $v1 = $range1.Value2
    $v2 = $range2.Value2

    $range1.Value2 = $v2
    $range2.Value2 = $v1
}
